# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.096.18"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").Value = "1.657.03"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").Value = "'218.53"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6
$ws.Range("D6").Value = "'0.5301"
$ws.Range("E6").Value = "  +1.66%  "

# Row 7
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").Value = "'0.2617"
$ws.Range("E8").Value = "  -1.92%  "

# Row 9
$ws.Range("D9").Value = "'0.06339"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").Value = "'20.44"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11
$ws.Range("D11").Value = "'0.07771"
$ws.Range("E11").Value = "  +0.85%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.677.90"
$ws.Range("E12").Value = "  +0.71%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.501"
$ws.Range("E13").Value = "  +1.70%  "

# Row 14
$ws.Range("D14").Value = "'0.5479"
$ws.Range("E14").Value = "  +0.19%  "

# Row 15
$ws.Range("D15").Value = "0.0₅8169"
$ws.Range("E15").Value = "  -0.45%  "

# Row 16
$ws.Range("D16").Value = "'65.28"

# Row 17
$ws.Range("D17").Value = "26.125.48"
$ws.Range("E17").Value = "  -0.42%  "

# Row 18
$ws.Range("E18").Value = "  -0.32%  "

# Row 19
$ws.Range("D19").Value = "'4.549"
$ws.Range("E19").Value = "  -2.07%  "

# Row 20
$ws.Range("D20").Value = "'193.61"
$ws.Range("E20").Value = "  -0.73%  "

# Row 21
$ws.Range("E21").Value = "  -0.44%  "

# Row 22
$ws.Range("D22").Value = "'6.022"
$ws.Range("E22").Value = "  -0.99%  "

# Row 23
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  -0.42%  "

# Row 24
$ws.Range("D24").Value = "'140.41"
$ws.Range("E24").Value = "  +0.90%  "

# Row 25
$ws.Range("D25").Value = "'0.1246"
$ws.Range("E25").Value = "  +0.37%  "

# Row 26
$ws.Range("D26").Value = "'7.282"
$ws.Range("E26").Value = "  +0.78%  "

# Row 27
$ws.Range("D27").Value = "'16.18"
$ws.Range("E27").Value = "  -0.06%  "

# Row 29
$ws.Range("D29").Value = "'0.05951"
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("D30").Value = "'1.277"
$ws.Range("E30").Value = "  -0.37%  "

# Row 31
$ws.Range("D31").Value = "'3.512"
$ws.Range("E31").Value = "  -3.26%  "

# Row 32
$ws.Range("D32").Value = "'3.243"
$ws.Range("E32").Value = "  -1.84%  "

# Row 33
$ws.Range("D33").Value = "'1.557"
$ws.Range("E33").Value = "  -4.50%  "

# Row 34
$ws.Range("D34").Value = "'0.9493"
$ws.Range("E34").Value = "  -2.94%  "

# Row 35
$ws.Range("D35").Value = "'2.412"
$ws.Range("E35").Value = "  -0.46%  "

# Row 36
$ws.Range("D36").Value = "'2.771"
$ws.Range("E36").Value = "  -0.25%  "

# Row 37
$ws.Range("D37").Value = "'0.5642"

# Row 38
$ws.Range("D38").Value = "'0.01613"
$ws.Range("E38").Value = "  +1.11%  "

# Row 39
$ws.Range("E39").Value = "  -2.13%  "

# Row 40
$ws.Range("D40").Value = "'0.8481"
$ws.Range("E40").Value = "  -0.97%  "

# Row 41
$ws.Range("E41").Value = "  -0.26%  "

# Row 42
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'101.44"
$ws.Range("E42").Value = "  +1.56%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.013.26"
$ws.Range("E43").Value = "  -1.26%  "

# Row 44
$ws.Range("D44").Value = "1.801.52"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45
$ws.Range("D45").Value = "'57.11"
$ws.Range("E45").Value = "  -0.23%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'1.004"
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈104"
$ws.Range("E47").Value = "  -5.04%  "

# Row 48
$ws.Range("D48").Value = "'0.4286"
$ws.Range("E48").Value = "  +1.35%  "

# Row 49
$ws.Range("D49").Value = "'0.05155"
$ws.Range("E49").Value = "  -0.67%  "

# Row 50
$ws.Range("D50").Value = "'1.469"
$ws.Range("E50").Value = "  +0.31%  "

# Row 51
$ws.Range("D51").Value = "'7.722"
$ws.Range("E51").Value = "  -4.09%  "

